# Generate Report for Handoff
# Updates the Priority column (E) and Latest Handoff Datetime column (H)
# for the rows that were just handed off (rows 4-7) in both the zh-cn and
# de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7, Priority -> "ht", Latest Handoff Datetime -> 2016-09-06 16:57:23
foreach ($r in 4..7) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-06 16:57:23"
}

# de-de sheet: rows 4-7, Priority -> "ht", Latest Handoff Datetime -> 2016-09-06 16:57:30
foreach ($r in 4..7) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-06 16:57:30"
}

# The Overview sheet's "Latest HO Xliff Generate Date" column (G) shares the
# same underlying text as de-de's "Latest Handoff Datetime" for these rows,
# so it moves in lockstep with the de-de update above.
foreach ($r in 4..7) {
    $wsOverview.Range("G$r").Value = "2016-09-06 16:57:30"
}
